# Nayem_meal.xlsx - "Moyla correction & March 4&5"
# Fill in meal entries for the 4th and 5th (columns E and F) for each
# person, add a bazar entry on the 4th (E43), and switch the running
# "cost/deposit" rows (23-29) over to a live formula instead of a
# hard-coded constant.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Meal entries for March 4 (col E) and March 5 (col F), rows 3-9 ---
$ws.Range("E3").Value = 2.5
$ws.Range("F3").Value = 1.5

$ws.Range("E4").Value = 2.5
$ws.Range("F4").Value = 1.5

$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0

$ws.Range("E6").Value = 2.5
$ws.Range("F6").Value = 2.5

$ws.Range("E7").Value = 2.5
$ws.Range("F7").Value = 2.5

$ws.Range("E8").Value = 2.5
$ws.Range("F8").Value = 2.5

$ws.Range("E9").Value = 2.5
$ws.Range("F9").Value = 2.5

# --- Bazar entry added on the 4th ---
$ws.Range("E43").Value = 515

# --- cost/deposit per-day rows: replace hard-coded -14.285714 with a
#     live formula, and correct F23 ---
$ws.Range("E23").Formula = "=-100/7"
$ws.Range("E24:E29").Formula = "=-100/7"
$ws.Range("F23").Value = 15

# --- View state: scroll position & active selection ---
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E29").Select()
